$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw Data")
$ws.Activate()

# New pricing for Drift: update Tier 1-3 price/probability for Product 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.25

$ws.Range("D3").Value = 360
$ws.Range("E3").Value = 0.5

$ws.Range("D4").Value = 1200
$ws.Range("E4").Value = 0.25

# Remove the old Tier 4 row for Product 1 (row 5); rows below shift up
$ws.Rows.Item(5).Delete()

# Update selection to match new layout
$ws.Range("A5:XFD5").Select() | Out-Null
